# The sheet is a weekly price log for "Naranja" (orange) at the
# "Vega Monumental Concepción" market. The edit inserts one new daily
# record as row 94 (pushing the previous rows 94-171 down to 95-172),
# growing the used range from A1:T171 to A1:T172.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 94; existing rows 94-171 shift to 95-172.
$ws.Rows.Item(94).Insert()

# Populate the new row with the new record's data.
$ws.Range("A94").Value = 11
$ws.Range("B94").Value = "Vega Monumental Concepción"
$ws.Range("C94").Value = "Bíobío"
$ws.Range("D94").Value = 44512
$ws.Range("E94").Value = 8
$ws.Range("F94").Value = "Fruta"
$ws.Range("G94").Value = 100102
$ws.Range("H94").Value = "Cítricos"
$ws.Range("I94").Value = 100102005
$ws.Range("J94").Value = "Naranja"
$ws.Range("K94").Value = "Navel Late"
$ws.Range("L94").Value = "Primera"
$ws.Range("M94").Value = 340
$ws.Range("N94").Value = 9000
$ws.Range("O94").Value = 10000
$ws.Range("P94").Value = 9588
$ws.Range("Q94").Value = "$/bandeja 15 kilos granel"
$ws.Range("R94").Value = "Región de O'Higgins"
$ws.Range("S94").Value = 639
$ws.Range("T94").Value = 15
